$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(28, 8).Value = 1597  # H28: 1731.7 -> 1597
$ws.Cells.Item(28, 9).Value = 756.8  # I28: 813.1111 -> 756.8
$ws.Cells.Item(28, 11).Value = 756.8  # K28: 813.1111 -> 756.8
$ws.Cells.Item(28, 13).Value = -271.8  # M28: -328.1111 -> -271.8
$ws.Cells.Item(62, 8).Value = 7282.3335  # H62: 6990.125 -> 7282.3335
$ws.Cells.Item(62, 9).Value = 6897.467  # I62: 6340.5 -> 6897.467
$ws.Cells.Item(62, 10).Value = 7923.778  # J62: 8289.375 -> 7923.778
$ws.Cells.Item(62, 11).Value = 6897.467  # K62: 6340.5 -> 6897.467
$ws.Cells.Item(62, 12).Value = 7923.778  # L62: 8289.375 -> 7923.778
$ws.Cells.Item(62, 13).Value = -6273.467  # M62: -5716.5 -> -6273.467
$ws.Cells.Item(62, 14).Value = -9171.778  # N62: -9537.375 -> -9171.778
$ws.Cells.Item(65, 8).Value = 7282.3335  # H65: 6990.125 -> 7282.3335
$ws.Cells.Item(65, 9).Value = 6897.467  # I65: 6340.5 -> 6897.467
$ws.Cells.Item(65, 10).Value = 7923.778  # J65: 8289.375 -> 7923.778
$ws.Cells.Item(65, 11).Value = 34487.335  # K65: 31702.5 -> 34487.335
$ws.Cells.Item(65, 12).Value = 39618.89  # L65: 41446.875 -> 39618.89
$ws.Cells.Item(65, 13).Value = -31367.335  # M65: -28582.5 -> -31367.335
$ws.Cells.Item(65, 14).Value = -45858.89  # N65: -47686.875 -> -45858.89
$ws.Cells.Item(70, 8).Value = 2930.2222  # H70: 4722.852 -> 2930.2222
$ws.Cells.Item(70, 9).Value = 1914.3334  # I70: 1915.4445 -> 1914.3334
$ws.Cells.Item(70, 10).Value = 3438.1667  # J70: 6126.5557 -> 3438.1667
$ws.Cells.Item(70, 11).Value = 5743.0002  # K70: 5746.333500000001 -> 5743.0002
$ws.Cells.Item(70, 12).Value = 10314.5001  # L70: 18379.6671 -> 10314.5001
$ws.Cells.Item(70, 13).Value = -5473.0002  # M70: -5476.333500000001 -> -5473.0002
$ws.Cells.Item(70, 14).Value = -10854.5001  # N70: -18919.6671 -> -10854.5001
$ws.Cells.Item(73, 8).Value = 2930.2222  # H73: 4722.852 -> 2930.2222
$ws.Cells.Item(73, 9).Value = 1914.3334  # I73: 1915.4445 -> 1914.3334
$ws.Cells.Item(73, 10).Value = 3438.1667  # J73: 6126.5557 -> 3438.1667
$ws.Cells.Item(73, 11).Value = 5743.0002  # K73: 5746.333500000001 -> 5743.0002
$ws.Cells.Item(73, 12).Value = 10314.5001  # L73: 18379.6671 -> 10314.5001
$ws.Cells.Item(73, 13).Value = -4807.0002  # M73: -4810.333500000001 -> -4807.0002
$ws.Cells.Item(73, 14).Value = -12186.5001  # N73: -20251.6671 -> -12186.5001
$ws.Cells.Item(76, 8).Value = 4705.162  # H76: 4747.081 -> 4705.162
$ws.Cells.Item(76, 9).Value = 4111.905  # I76: 4235.05 -> 4111.905
$ws.Cells.Item(76, 10).Value = 5483.8125  # J76: 5349.4707 -> 5483.8125
$ws.Cells.Item(76, 11).Value = 4111.905  # K76: 4235.05 -> 4111.905
$ws.Cells.Item(76, 12).Value = 5483.8125  # L76: 5349.4707 -> 5483.8125
$ws.Cells.Item(76, 13).Value = -3796.905  # M76: -3920.05 -> -3796.905
$ws.Cells.Item(76, 14).Value = -6113.8125  # N76: -5979.4707 -> -6113.8125
$ws.Cells.Item(79, 8).Value = 4705.162  # H79: 4747.081 -> 4705.162
$ws.Cells.Item(79, 9).Value = 4111.905  # I79: 4235.05 -> 4111.905
$ws.Cells.Item(79, 10).Value = 5483.8125  # J79: 5349.4707 -> 5483.8125
$ws.Cells.Item(79, 11).Value = 4111.905  # K79: 4235.05 -> 4111.905
$ws.Cells.Item(79, 12).Value = 5483.8125  # L79: 5349.4707 -> 5483.8125
$ws.Cells.Item(79, 13).Value = -3019.905  # M79: -3143.05 -> -3019.905
$ws.Cells.Item(79, 14).Value = -7667.8125  # N79: -7533.4707 -> -7667.8125
$ws.Cells.Item(98, 8).Value = 1491.8077  # H98: 1884.4375 -> 1491.8077
$ws.Cells.Item(98, 9).Value = 994.1818  # I98: 1048.9231 -> 994.1818
$ws.Cells.Item(98, 10).Value = 4228.75  # J98: 5505 -> 4228.75
$ws.Cells.Item(98, 11).Value = 994.1818  # K98: 1048.9231 -> 994.1818
$ws.Cells.Item(98, 12).Value = 4228.75  # L98: 5505 -> 4228.75
$ws.Cells.Item(98, 13).Value = 503.8182  # M98: 449.0769 -> 503.8182
$ws.Cells.Item(98, 14).Value = -7224.75  # N98: -8501 -> -7224.75
$ws.Cells.Item(106, 8).Value = 3091.6667  # H106: 4466.5 -> 3091.6667
$ws.Cells.Item(106, 9).Value = 2442.8572  # I106: 2100 -> 2442.8572
$ws.Cells.Item(106, 10).Value = 4000  # J106: 4939.8 -> 4000
$ws.Cells.Item(106, 11).Value = 2442.8572  # K106: 2100 -> 2442.8572
$ws.Cells.Item(106, 12).Value = 4000  # L106: 4939.8 -> 4000
$ws.Cells.Item(106, 13).Value = -1811.8572  # M106: -1469 -> -1811.8572
$ws.Cells.Item(106, 14).Value = -5262  # N106: -6201.8 -> -5262
$ws.Cells.Item(107, 8).Value = 1786.2273  # H107: 1867.0952 -> 1786.2273
$ws.Cells.Item(107, 9).Value = 1564.85  # I107: 1642.579 -> 1564.85
$ws.Cells.Item(107, 11).Value = 1564.85  # K107: 1642.579 -> 1564.85
$ws.Cells.Item(107, 13).Value = 355.1500000000001  # M107: 277.421 -> 355.1500000000001
$ws.Cells.Item(116, 8).Value = 5121.9287  # H116: 5431.3076 -> 5121.9287
$ws.Cells.Item(116, 9).Value = 3708.5  # I116: 4081.1428 -> 3708.5
$ws.Cells.Item(116, 11).Value = 3708.5  # K116: 4081.1428 -> 3708.5
$ws.Cells.Item(116, 13).Value = -266.5  # M116: -639.1428000000001 -> -266.5
$ws.Cells.Item(122, 8).Value = 1491.8077  # H122: 1884.4375 -> 1491.8077
$ws.Cells.Item(122, 9).Value = 994.1818  # I122: 1048.9231 -> 994.1818
$ws.Cells.Item(122, 10).Value = 4228.75  # J122: 5505 -> 4228.75
$ws.Cells.Item(122, 11).Value = 2982.5454  # K122: 3146.7693 -> 2982.5454
$ws.Cells.Item(122, 12).Value = 12686.25  # L122: 16515 -> 12686.25
$ws.Cells.Item(122, 13).Value = -532.5454  # M122: -696.7692999999999 -> -532.5454
$ws.Cells.Item(122, 14).Value = -17586.25  # N122: -21415 -> -17586.25
$ws.Cells.Item(132, 8).Value = 2992160.8  # H132: 3191621.8 -> 2992160.8
$ws.Cells.Item(132, 9).Value = 3191510  # I132: 3339941.2 -> 3191510
$ws.Cells.Item(132, 10).Value = 1921.6666  # J132: 2749.5 -> 1921.6666
$ws.Cells.Item(132, 11).Value = 9574530  # K132: 10019823.6 -> 9574530
$ws.Cells.Item(132, 12).Value = 5764.9998  # L132: 8248.5 -> 5764.9998
$ws.Cells.Item(132, 13).Value = -9572000  # M132: -10017293.6 -> -9572000
$ws.Cells.Item(132, 14).Value = -10824.9998  # N132: -13308.5 -> -10824.9998
$ws.Cells.Item(137, 8).Value = 9286.441000000001  # H137: 8756.695 -> 9286.441000000001
$ws.Cells.Item(137, 9).Value = 12845.777  # I137: 11677.566 -> 12845.777
$ws.Cells.Item(137, 11).Value = 38537.331  # K137: 35032.698 -> 38537.331
$ws.Cells.Item(137, 13).Value = -35987.331  # M137: -32482.698 -> -35987.331
$ws.Cells.Item(138, 8).Value = 2804.2856  # H138: 2871.1482 -> 2804.2856
$ws.Cells.Item(138, 9).Value = 1914.3125  # I138: 1975.3334 -> 1914.3125
$ws.Cells.Item(138, 11).Value = 5742.9375  # K138: 5926.0002 -> 5742.9375
$ws.Cells.Item(138, 13).Value = -602.9375  # M138: -786.0002000000004 -> -602.9375
$ws.Cells.Item(141, 8).Value = 1666  # H141: 1799.4 -> 1666
$ws.Cells.Item(141, 9).Value = 1332.6666  # I141: 1499.5 -> 1332.6666
$ws.Cells.Item(141, 11).Value = 3997.9998  # K141: 4498.5 -> 3997.9998
$ws.Cells.Item(141, 13).Value = 1182.0002  # M141: 681.5 -> 1182.0002

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(61, 8).Value = 4748.5  # H61: 4642.4443 -> 4748.5
$ws.Cells.Item(61, 9).Value = 1997.5  # I61: 2356.8 -> 1997.5
$ws.Cells.Item(61, 11).Value = 1997.5  # K61: 2356.8 -> 1997.5
$ws.Cells.Item(61, 13).Value = -1785.5  # M61: -2144.8 -> -1785.5
$ws.Cells.Item(110, 8).Value = 1167.3846  # H110: 1267 -> 1167.3846
$ws.Cells.Item(110, 9).Value = 1139.6666  # I110: 1243.7 -> 1139.6666
$ws.Cells.Item(110, 11).Value = 1139.6666  # K110: 1243.7 -> 1139.6666
$ws.Cells.Item(110, 13).Value = 905.3334  # M110: 801.3 -> 905.3334
$ws.Cells.Item(122, 8).Value = 2349.4443  # H122: 1854.9231 -> 2349.4443
$ws.Cells.Item(122, 9).Value = 2018.125  # I122: 1592.8334 -> 2018.125
$ws.Cells.Item(122, 11).Value = 6054.375  # K122: 4778.5002 -> 6054.375
$ws.Cells.Item(122, 13).Value = -3604.375  # M122: -2328.5002 -> -3604.375
$ws.Cells.Item(132, 8).Value = 3933  # H132: 3798.3635 -> 3933
$ws.Cells.Item(132, 9).Value = 4403  # I132: 4032.75 -> 4403
$ws.Cells.Item(132, 11).Value = 13209  # K132: 12098.25 -> 13209
$ws.Cells.Item(132, 13).Value = -10679  # M132: -9568.25 -> -10679
$ws.Cells.Item(136, 8).Value = 4748.5  # H136: 4642.4443 -> 4748.5
$ws.Cells.Item(136, 9).Value = 1997.5  # I136: 2356.8 -> 1997.5
$ws.Cells.Item(136, 11).Value = 5992.5  # K136: 7070.400000000001 -> 5992.5
$ws.Cells.Item(136, 13).Value = -3442.5  # M136: -4520.400000000001 -> -3442.5
$ws.Cells.Item(140, 8).Value = 67231  # H140: 67533 -> 67231
$ws.Cells.Item(140, 10).Value = 67231  # J140: 67533 -> 67231
$ws.Cells.Item(140, 12).Value = 67231  # L140: 67533 -> 67231
$ws.Cells.Item(140, 14).Value = -77591  # N140: -77893 -> -77591

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(99, 8).Value = 2602.1667  # H99: 2643.125 -> 2602.1667
$ws.Cells.Item(99, 9).Value = 1393.3  # I99: 1488.2222 -> 1393.3
$ws.Cells.Item(99, 10).Value = 4113.25  # J99: 4128 -> 4113.25
$ws.Cells.Item(99, 11).Value = 1393.3  # K99: 1488.2222 -> 1393.3
$ws.Cells.Item(99, 12).Value = 4113.25  # L99: 4128 -> 4113.25
$ws.Cells.Item(99, 13).Value = 104.7  # M99: 9.77780000000007 -> 104.7
$ws.Cells.Item(99, 14).Value = -7109.25  # N99: -7124 -> -7109.25
$ws.Cells.Item(134, 8).Value = 4236.2856  # H134: 3815.6875 -> 4236.2856
$ws.Cells.Item(134, 9).Value = 4236.2856  # I134: 3815.6875 -> 4236.2856
$ws.Cells.Item(134, 11).Value = 12708.8568  # K134: 11447.0625 -> 12708.8568
$ws.Cells.Item(134, 13).Value = -10173.8568  # M134: -8912.0625 -> -10173.8568

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(16, 8).Value = 2385.0715  # H16: 2406.1538 -> 2385.0715
$ws.Cells.Item(16, 10).Value = 3044.2  # J16: 3277.5 -> 3044.2
$ws.Cells.Item(16, 12).Value = 3044.2  # L16: 3277.5 -> 3044.2
$ws.Cells.Item(16, 14).Value = -3618.2  # N16: -3851.5 -> -3618.2
$ws.Cells.Item(31, 8).Value = 2130486.8  # H31: 1925603 -> 2130486.8
$ws.Cells.Item(31, 9).Value = 2566077.5  # I31: 2223976.2 -> 2566077.5
$ws.Cells.Item(31, 10).Value = 6982.375  # J31: 7489.2856 -> 6982.375
$ws.Cells.Item(31, 11).Value = 2566077.5  # K31: 2223976.2 -> 2566077.5
$ws.Cells.Item(31, 12).Value = 6982.375  # L31: 7489.2856 -> 6982.375
$ws.Cells.Item(31, 13).Value = -2565782.5  # M31: -2223681.2 -> -2565782.5
$ws.Cells.Item(31, 14).Value = -7572.375  # N31: -8079.2856 -> -7572.375
$ws.Cells.Item(34, 8).Value = 2130486.8  # H34: 1925603 -> 2130486.8
$ws.Cells.Item(34, 9).Value = 2566077.5  # I34: 2223976.2 -> 2566077.5
$ws.Cells.Item(34, 10).Value = 6982.375  # J34: 7489.2856 -> 6982.375
$ws.Cells.Item(34, 11).Value = 2566077.5  # K34: 2223976.2 -> 2566077.5
$ws.Cells.Item(34, 12).Value = 6982.375  # L34: 7489.2856 -> 6982.375
$ws.Cells.Item(34, 13).Value = -2565875.5  # M34: -2223774.2 -> -2565875.5
$ws.Cells.Item(34, 14).Value = -7386.375  # N34: -7893.2856 -> -7386.375
$ws.Cells.Item(103, 8).Value = 5309.4  # H103: 9444.333000000001 -> 5309.4
$ws.Cells.Item(103, 9).Value = 5309.4  # I103: 6256 -> 5309.4
$ws.Cells.Item(103, 10).Value = 0  # J103: 11995 -> 0
$ws.Cells.Item(103, 11).Value = 5309.4  # K103: 6256 -> 5309.4
$ws.Cells.Item(103, 12).Value = 0  # L103: 11995 -> 0
$ws.Cells.Item(103, 13).Value = -4137.4  # M103: -5084 -> -4137.4
$ws.Cells.Item(103, 14).ClearContents()  # N103: delete (was -14339)
$ws.Cells.Item(105, 8).Value = 2300  # H105: 3006.5715 -> 2300
$ws.Cells.Item(105, 10).Value = 1392.8572  # J105: 1949 -> 1392.8572
$ws.Cells.Item(105, 12).Value = 1392.8572  # L105: 1949 -> 1392.8572
$ws.Cells.Item(105, 14).Value = -4886.8572  # N105: -5443 -> -4886.8572
$ws.Cells.Item(107, 8).Value = 411.14285  # H107: 439.07693 -> 411.14285
$ws.Cells.Item(107, 9).Value = 284  # I107: 313.5 -> 284
$ws.Cells.Item(107, 11).Value = 284  # K107: 313.5 -> 284
$ws.Cells.Item(107, 13).Value = 1636  # M107: 1606.5 -> 1636
$ws.Cells.Item(113, 8).Value = 2385.0715  # H113: 2406.1538 -> 2385.0715
$ws.Cells.Item(113, 10).Value = 3044.2  # J113: 3277.5 -> 3044.2
$ws.Cells.Item(113, 12).Value = 3044.2  # L113: 3277.5 -> 3044.2
$ws.Cells.Item(113, 14).Value = -7384.2  # N113: -7617.5 -> -7384.2
$ws.Cells.Item(122, 8).Value = 8151.6333  # H122: 8422.379000000001 -> 8151.6333
$ws.Cells.Item(122, 9).Value = 9527.708000000001  # I122: 9928.913 -> 9527.708000000001
$ws.Cells.Item(122, 11).Value = 28583.124  # K122: 29786.739 -> 28583.124
$ws.Cells.Item(122, 13).Value = -26133.124  # M122: -27336.739 -> -26133.124
$ws.Cells.Item(134, 8).Value = 3820.8  # H134: 2373.6365 -> 3820.8
$ws.Cells.Item(134, 9).Value = 3820.8  # I134: 2061 -> 3820.8
$ws.Cells.Item(134, 10).Value = 0  # J134: 5500 -> 0
$ws.Cells.Item(134, 11).Value = 11462.4  # K134: 6183 -> 11462.4
$ws.Cells.Item(134, 12).Value = 0  # L134: 16500 -> 0
$ws.Cells.Item(134, 13).Value = -8927.400000000001  # M134: -3648 -> -8927.400000000001
$ws.Cells.Item(134, 14).ClearContents()  # N134: delete (was -21570)

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(37, 8).Value = 39953.47  # H37: 39953.277 -> 39953.47
$ws.Cells.Item(37, 10).Value = 39953.47  # J37: 39953.277 -> 39953.47
$ws.Cells.Item(37, 12).Value = 119860.41  # L37: 119859.831 -> 119860.41
$ws.Cells.Item(37, 14).Value = -120084.41  # N37: -120083.831 -> -120084.41
$ws.Cells.Item(88, 8).Value = 7776  # H88: 5925.3335 -> 7776
$ws.Cells.Item(88, 10).Value = 7775  # J88: 5396.4287 -> 7775
$ws.Cells.Item(88, 12).Value = 23325  # L88: 16189.2861 -> 23325
$ws.Cells.Item(88, 14).Value = -24181  # N88: -17045.2861 -> -24181
$ws.Cells.Item(91, 8).Value = 7776  # H91: 5925.3335 -> 7776
$ws.Cells.Item(91, 10).Value = 7775  # J91: 5396.4287 -> 7775
$ws.Cells.Item(91, 12).Value = 23325  # L91: 16189.2861 -> 23325
$ws.Cells.Item(91, 14).Value = -26289  # N91: -19153.2861 -> -26289
$ws.Cells.Item(139, 8).Value = 10134.777  # H139: 11252.875 -> 10134.777
$ws.Cells.Item(139, 9).Value = 10147.5  # I139: 11427.143 -> 10147.5
$ws.Cells.Item(139, 11).Value = 30442.5  # K139: 34281.429 -> 30442.5
$ws.Cells.Item(139, 13).Value = -25302.5  # M139: -29141.429 -> -25302.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(139, 8).Value = 80331.336  # H139: 71498.5 -> 80331.336
$ws.Cells.Item(139, 10).Value = 80331.336  # J139: 71498.5 -> 80331.336
$ws.Cells.Item(139, 12).Value = 80331.336  # L139: 71498.5 -> 80331.336
$ws.Cells.Item(139, 14).Value = -90611.336  # N139: -81778.5 -> -90611.336
$ws.Cells.Item(140, 8).Value = 0  # H140: 189750 -> 0
$ws.Cells.Item(140, 10).Value = 0  # J140: 189750 -> 0
$ws.Cells.Item(140, 12).Value = 0  # L140: 189750 -> 0
$ws.Cells.Item(140, 14).ClearContents()  # N140: delete (was -200110)

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(45, 8).Value = 48046  # H45: 36682 -> 48046
$ws.Cells.Item(45, 9).Value = 0  # I45: 14000 -> 0
$ws.Cells.Item(45, 10).Value = 48046  # J45: 48023 -> 48046
$ws.Cells.Item(45, 11).Value = 0  # K45: 14000 -> 0
$ws.Cells.Item(45, 12).Value = 48046  # L45: 48023 -> 48046
$ws.Cells.Item(45, 13).ClearContents()  # M45: delete (was -13593)
$ws.Cells.Item(45, 14).Value = -48860  # N45: -48837 -> -48860
$ws.Cells.Item(61, 8).Value = 1297  # H61: 1246.25 -> 1297
$ws.Cells.Item(61, 9).Value = 1150.3846  # I61: 1102.8572 -> 1150.3846
$ws.Cells.Item(61, 11).Value = 1150.3846  # K61: 1102.8572 -> 1150.3846
$ws.Cells.Item(61, 13).Value = -948.3846000000001  # M61: -900.8571999999999 -> -948.3846000000001
$ws.Cells.Item(82, 8).Value = 10752.962  # H82: 11320.75 -> 10752.962
$ws.Cells.Item(82, 9).Value = 9873  # I82: 10720.643 -> 9873
$ws.Cells.Item(82, 11).Value = 9873  # K82: 10720.643 -> 9873
$ws.Cells.Item(82, 13).Value = -9512  # M82: -10359.643 -> -9512
$ws.Cells.Item(85, 8).Value = 10752.962  # H85: 11320.75 -> 10752.962
$ws.Cells.Item(85, 9).Value = 9873  # I85: 10720.643 -> 9873
$ws.Cells.Item(85, 11).Value = 9873  # K85: 10720.643 -> 9873
$ws.Cells.Item(85, 13).Value = -8625  # M85: -9472.643 -> -8625
$ws.Cells.Item(113, 8).Value = 1297  # H113: 1246.25 -> 1297
$ws.Cells.Item(113, 9).Value = 1150.3846  # I113: 1102.8572 -> 1150.3846
$ws.Cells.Item(113, 11).Value = 1150.3846  # K113: 1102.8572 -> 1150.3846
$ws.Cells.Item(113, 13).Value = 1019.6154  # M113: 1067.1428 -> 1019.6154
$ws.Cells.Item(123, 8).Value = 78999  # H123: 79494 -> 78999
$ws.Cells.Item(123, 10).Value = 78999  # J123: 79494 -> 78999
$ws.Cells.Item(123, 12).Value = 78999  # L123: 79494 -> 78999
$ws.Cells.Item(123, 14).Value = -88799  # N123: -89294 -> -88799

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(74, 8).Value = 16209.286  # H74: 14486.125 -> 16209.286
$ws.Cells.Item(74, 9).Value = 10569  # I74: 8711.333000000001 -> 10569
$ws.Cells.Item(74, 10).Value = 17149.334  # J74: 17951 -> 17149.334
$ws.Cells.Item(74, 11).Value = 10569  # K74: 8711.333000000001 -> 10569
$ws.Cells.Item(74, 12).Value = 17149.334  # L74: 17951 -> 17149.334
$ws.Cells.Item(74, 13).Value = -9633  # M74: -7775.333000000001 -> -9633
$ws.Cells.Item(74, 14).Value = -19021.334  # N74: -19823 -> -19021.334
$ws.Cells.Item(77, 8).Value = 16209.286  # H77: 14486.125 -> 16209.286
$ws.Cells.Item(77, 9).Value = 10569  # I77: 8711.333000000001 -> 10569
$ws.Cells.Item(77, 10).Value = 17149.334  # J77: 17951 -> 17149.334
$ws.Cells.Item(77, 11).Value = 31707  # K77: 26133.999 -> 31707
$ws.Cells.Item(77, 12).Value = 51448.00199999999  # L77: 53853 -> 51448.00199999999
$ws.Cells.Item(77, 13).Value = -27027  # M77: -21453.999 -> -27027
$ws.Cells.Item(77, 14).Value = -60808.00199999999  # N77: -63213 -> -60808.00199999999
$ws.Cells.Item(122, 8).Value = 94762.24000000001  # H122: 66146.44500000001 -> 94762.24000000001
$ws.Cells.Item(122, 9).Value = 137244.83  # I122: 86811.03999999999 -> 137244.83
$ws.Cells.Item(122, 10).Value = 4486.75  # J122: 4152.6665 -> 4486.75
$ws.Cells.Item(122, 11).Value = 411734.49  # K122: 260433.12 -> 411734.49
$ws.Cells.Item(122, 12).Value = 13460.25  # L122: 12457.9995 -> 13460.25
$ws.Cells.Item(122, 13).Value = -409284.49  # M122: -257983.12 -> -409284.49
$ws.Cells.Item(122, 14).Value = -18360.25  # N122: -17357.9995 -> -18360.25
